$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values to match repulled/pushed data and recalculated mean
$ws.Range("F9").Value  = -8
$ws.Range("F12").Value = 7
$ws.Range("F14").Value = -3
$ws.Range("F15").Value = 11
$ws.Range("F17").Value = 0
$ws.Range("F20").Value = -7
$ws.Range("F27").Value = 0
$ws.Range("F28").Value = -5
$ws.Range("F29").Value = -2
$ws.Range("F31").Value = -1
$ws.Range("F33").Value = -4
